# Update gh-pages output (漫展信息/con-expo tracker): bump the "want to go"
# counts (column F) on a handful of rows across the four sheets, matching a
# fresh scrape/regeneration of the data (commit "Update gh-pages to output
# generated at 456a3b4"). No other cell content changes.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 54661
$ws1.Range("F12").Value = 5210
$ws1.Range("F19").Value = 1264
$ws1.Range("F23").Value = 359
$ws1.Range("F29").Value = 4958
$ws1.Range("F31").Value = 4919
$ws1.Range("F32").Value = 8911
$ws1.Range("F33").Value = 113
$ws1.Range("F40").Value = 4194

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 1127

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 784

# --- Sheet "全部类型" (All types, aggregated view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 784
$ws4.Range("F21").Value = 1264
$ws4.Range("F26").Value = 359
$ws4.Range("F31").Value = 4958
$ws4.Range("F33").Value = 4919
$ws4.Range("F34").Value = 8911
$ws4.Range("F35").Value = 113
$ws4.Range("F44").Value = 4194
